$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column AI (18-jul) ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the header style from AH1 (bold/centered) onto the new AI1 header cell
$wsSpot.Range("AH1").Copy()
$wsSpot.Range("AI1").PasteSpecial(-4122)
$wsSpot.Range("AI1").Value = "18-jul"

$spotValues = @{
    2  = 100.97
    3  = 92.22
    4  = 91.20999999999999
    5  = 82.52
    6  = 78.45
    7  = 84.31
    8  = 90.31999999999999
    9  = 105.08
    10 = 109.39
    11 = 90.64
    12 = 80
    13 = 81.58
    14 = 80
    15 = 73.02
    16 = 71.2
    17 = 72.68000000000001
    18 = 79.40000000000001
    19 = 83.22
    20 = 93.97
    21 = 108.1
    22 = 110
    23 = 119.39
    24 = 119.19
    25 = 111.65
}

foreach ($row in $spotValues.Keys) {
    $wsSpot.Cells.Item($row, 35).Value = $spotValues[$row]
}

# --- Sheet "Gaz": add row 32 (2025-07-16, 34.8) ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A32").NumberFormat = "@"
$wsGaz.Range("A32").Value = "2025-07-16"
$wsGaz.Range("A32").Style = "Normal"
$wsGaz.Range("B32").Value = 34.8

# --- Sheet "CO2": add row 32 (2025-07-16, 70.39) ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A32").NumberFormat = "@"
$wsCo2.Range("A32").Value = "2025-07-16"
$wsCo2.Range("A32").Style = "Normal"
$wsCo2.Range("B32").Value = 70.39
